$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 304, shifting existing rows 304..343 down to 305..344
$ws.Rows.Item(304).Insert(-4121)

# Populate the newly inserted row 304 with this week's data point.
# Fixed / repeated columns (copied from the surrounding rows of this market/product block).
$ws.Range("A304").Value2 = 10
$ws.Range("B304").Value2 = "Vega Modelo de Temuco"
$ws.Range("C304").Value2 = "La Araucanía"
$ws.Range("D304").Value2 = 45154
$ws.Range("E304").Value2 = 9
$ws.Range("F304").Value2 = 100114007
$ws.Range("G304").Value2 = "Jengibre"
$ws.Range("H304").Value2 = "Sin especificar"
$ws.Range("I304").Value2 = "Primera"
$ws.Range("J304").Value2 = 90
$ws.Range("K304").Value2 = 22000
$ws.Range("L304").Value2 = 24000
$ws.Range("M304").Value2 = 23111
$ws.Range("N304").Value2 = "`$/caja 13 kilos"
$ws.Range("O304").Value2 = "Perú"
$ws.Range("P304").Value2 = 1778
$ws.Range("Q304").Value2 = 13
$ws.Range("R304").Value2 = "Hortaliza"
